$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Re-format the I column (kj_min) to match the D/E/F "0.00" style ---
# Main body (style index 4 once saved): rows 2-18, 20-63, 65-70
$ws.Range("I2:I18").NumberFormat = "0.00"
$ws.Range("I20:I63").NumberFormat = "0.00"
$ws.Range("I65:I70").NumberFormat = "0.00"

# The two "NA" rows (19 and 64) keep right alignment like the other NA cells
# on those rows (style index 5 once saved)
$ws.Range("I19").NumberFormat = "0.00"
$ws.Range("I19").HorizontalAlignment = -4152
$ws.Range("I64").NumberFormat = "0.00"
$ws.Range("I64").HorizontalAlignment = -4152

# --- New K/L "statistikktabell" columns ---
# Populate K3:K13 and L3:L13 with the new data, formatted as whole numbers
$ws.Range("K3").Value = 1256.68
$ws.Range("L3").Value = 1324.56

$ws.Range("K4").Value = 1349.72
$ws.Range("L4").Value = 1198.72

$ws.Range("K5").Value = 1191.8800000000001
$ws.Range("L5").Value = 1095.52

$ws.Range("K6").Value = 1044.08
$ws.Range("L6").Value = 1117.44

$ws.Range("K7").Value = 1508.4
$ws.Range("L7").Value = 1503.76

$ws.Range("K8").Value = 1738.3157894736842
$ws.Range("L8").Value = 1377.2

$ws.Range("K9").Value = 1397.2
$ws.Range("L9").Value = 1592.72

$ws.Range("K10").Value = 1155.1600000000001
$ws.Range("L10").Value = 1305.96

$ws.Range("K11").Value = 1149.92
$ws.Range("L11").Value = 1314.84

$ws.Range("K12").Value = 1242.6400000000001
$ws.Range("L12").Value = 1175.92

$ws.Range("K13").Value = 1135.28
$ws.Range("L13").Value = 1145.56

$ws.Range("K3:L13").NumberFormat = "0"

# Two stray formatted-but-empty cells further down the (still unfinished) table
$ws.Range("K16").NumberFormat = "0.00"
$ws.Range("K18").NumberFormat = "0.00"

# --- View state: scrolled right one column, selection sitting on the new data ---
$ws.Range("L3:L13").Select() | Out-Null
$excel.ActiveWindow.ScrollColumn = 2

$wb.Save()
